{"js": "// Docx writer: omit jc attribute on table cells with AlignDefault.\n//\n// Every paragraph inside every table cell in this document carries an\n// explicit <w:jc w:val=\"left\"/> (AlignDefault). Left is the default\n// justification, so it should not be written at all. Where that\n// <w:jc> was the *only* child of <w:pPr>, removing it also means the\n// paragraph falls back onto whatever style the writer assigns instead\n// (the third table's cells had no pStyle at all, so the writer's\n// default paragraph styles - \"First Paragraph\" for the very first\n// paragraph of the table, \"Body Text\" for every other paragraph in it\n// - now show up explicitly in their place).\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet firstParagraphSeen = false;\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const rows = tables.items[t].rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (let r = 0; r < rows.items.length; r++) {\n    const cells = rows.items[r].cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (let c = 0; c < cells.items.length; c++) {\n      const paras = cells.items[c].body.paragraphs;\n      paras.load(\"items,style\");\n      await context.sync();\n\n      for (let p = 0; p < paras.items.length; p++) {\n        const para = paras.items[p];\n\n        // Dropping the redundant \"left\" alignment is all that the\n        // first two tables need - their paragraphs already carry an\n        // explicit pStyle (\"Compact\").\n        para.alignment = \"Left\";\n\n        // The third table's paragraphs have no pStyle at all, so the\n        // writer now falls back to its default paragraph styles in\n        // place of the dropped <w:jc>.\n        if (para.style === \"Normal\") {\n          para.style = firstParagraphSeen ? \"Body Text\" : \"First Paragraph\";\n          firstParagraphSeen = true;\n        }\n      }\n      await context.sync();\n    }\n  }\n}\n", "ps1": "# Docx writer: omit jc attribute on table cells with AlignDefault.\n#\n# Every paragraph inside every table cell in this document carries an\n# explicit left justification (AlignDefault). Left is the default\n# justification, so it should not be written to the OOXML at all.\n# Re-assigning wdAlignParagraphLeft (0) is a safe no-op for paragraphs\n# that are already left-aligned without an explicit <w:jc> (i.e. every\n# paragraph outside the three tables), and it drops the redundant\n# <w:jc w:val=\"left\"/> for the ones that do carry it.\n#\n# Where that <w:jc> used to be the *only* child of <w:pPr>, dropping it\n# also means the paragraph falls back onto whatever style the writer\n# assigns instead. The third table's cells have no pStyle at all (the\n# first two tables already use the \"Compact\" style), so the writer's\n# default paragraph styles now show up explicitly in its place: \"First\n# Paragraph\" for the very first paragraph of that table, \"Body Text\"\n# for every other paragraph inside it.\n\n$d = $word.ActiveDocument\n\n$firstParagraphAssigned = $false\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n\n    $para.Alignment = 0   # wdAlignParagraphLeft\n\n    # Range.Text includes the trailing paragraph mark (and, for the\n    # last paragraph in a table cell, the end-of-cell mark too) - trim\n    # those off before checking whether the paragraph actually holds\n    # any text.\n    $text = $para.Range.Text.TrimEnd(\"`r\", \"`a\")\n\n    if ($para.Style.NameLocal -eq \"Normal\" -and $text -ne \"\") {\n        if (-not $firstParagraphAssigned) {\n            $para.Style = \"First Paragraph\"\n            $firstParagraphAssigned = $true\n        } else {\n            $para.Style = \"Body Text\"\n        }\n    }\n}\n"}
